$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to hold a single text label ("Articulos") in A1. Replace it
# with numeric marker values across A1:C1 so B1/C1 can be referenced by
# the new named ranges below.
$ws.Range("A1").Value = 1
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 3

# Name the data columns (B1 = "Descripcion" header cell, C1 = "Precio" header cell).
$sheetName = $ws.Name
$wb.Names.Add("Descripcion", "=$sheetName!`$B`$1")
$wb.Names.Add("Precio", "=$sheetName!`$C`$1")

# Clear the stale selection range that used to be persisted in the sheet view.
$ws.Range("A1").Select()
